$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 5332.3335
$ws.Range("I62").Value = 4410.625
$ws.Range("J62").Value = 6069.7
$ws.Range("K62").Value = 4410.625
$ws.Range("L62").Value = 6069.7
$ws.Range("M62").Value = -3786.625
$ws.Range("N62").Value = -7317.7

$ws.Range("H65").Value = 5332.3335
$ws.Range("I65").Value = 4410.625
$ws.Range("J65").Value = 6069.7
$ws.Range("K65").Value = 22053.125
$ws.Range("L65").Value = 30348.5
$ws.Range("M65").Value = -18933.125
$ws.Range("N65").Value = -36588.5

$ws.Range("H92").Value = 43478696
$ws.Range("I92").Value = 66666948
$ws.Range("K92").Value = 66666948
$ws.Range("M92").Value = -66665700

$ws.Range("H137").Value = 1514.5927
$ws.Range("I137").Value = 1528.7368
$ws.Range("J137").Value = 1481
$ws.Range("K137").Value = 4586.2104
$ws.Range("L137").Value = 4443
$ws.Range("M137").Value = -2036.2104
$ws.Range("N137").Value = -9543

$ws.Range("H141").Value = 2171.3076
$ws.Range("I141").Value = 1457.7
$ws.Range("J141").Value = 4550
$ws.Range("K141").Value = 4373.1
$ws.Range("L141").Value = 13650
$ws.Range("M141").Value = 806.8999999999996
$ws.Range("N141").Value = -24010


# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1946.0857
$ws.Range("I61").Value = 1503.6451
$ws.Range("J61").Value = 5375
$ws.Range("K61").Value = 1503.6451
$ws.Range("L61").Value = 5375
$ws.Range("M61").Value = -1291.6451
$ws.Range("N61").Value = -5799

$ws.Range("H74").Value = 76925656
$ws.Range("I74").Value = 111112280
$ws.Range("K74").Value = 111112280
$ws.Range("M74").Value = -111111406

$ws.Range("H77").Value = 76925656
$ws.Range("I77").Value = 111112280
$ws.Range("K77").Value = 555561400
$ws.Range("M77").Value = -555557032

$ws.Range("H132").Value = 9378.531999999999
$ws.Range("I132").Value = 996.675
$ws.Range("J132").Value = 24618.273
$ws.Range("K132").Value = 2990.025
$ws.Range("L132").Value = 73854.819
$ws.Range("M132").Value = -460.0249999999996
$ws.Range("N132").Value = -78914.819

$ws.Range("H136").Value = 1946.0857
$ws.Range("I136").Value = 1503.6451
$ws.Range("J136").Value = 5375
$ws.Range("K136").Value = 4510.9353
$ws.Range("L136").Value = 16125
$ws.Range("M136").Value = -1960.9353
$ws.Range("N136").Value = -21225


# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4330.115
$ws.Range("I134").Value = 5074.7617
$ws.Range("K134").Value = 15224.2851
$ws.Range("M134").Value = -12689.2851


# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3052.2666
$ws.Range("I31").Value = 2733.3572
$ws.Range("J31").Value = 3196.2903
$ws.Range("K31").Value = 2733.3572
$ws.Range("L31").Value = 3196.2903
$ws.Range("M31").Value = -2438.3572
$ws.Range("N31").Value = -3786.2903

$ws.Range("H34").Value = 3052.2666
$ws.Range("I34").Value = 2733.3572
$ws.Range("J34").Value = 3196.2903
$ws.Range("K34").Value = 2733.3572
$ws.Range("L34").Value = 3196.2903
$ws.Range("M34").Value = -2531.3572
$ws.Range("N34").Value = -3600.2903

$ws.Range("H58").Value = 20273.482
$ws.Range("I58").Value = 1861.6
$ws.Range("J58").Value = 31104
$ws.Range("K58").Value = 1861.6
$ws.Range("L58").Value = 31104
$ws.Range("M58").Value = -1658.6
$ws.Range("N58").Value = -31510

$ws.Range("H99").Value = 33337432
$ws.Range("I99").Value = 3399.4443
$ws.Range("J99").Value = 83338480
$ws.Range("K99").Value = 3399.4443
$ws.Range("L99").Value = 83338480
$ws.Range("M99").Value = -1901.4443
$ws.Range("N99").Value = -83341476

$ws.Range("H126").Value = 33337432
$ws.Range("I126").Value = 3399.4443
$ws.Range("J126").Value = 83338480
$ws.Range("K126").Value = 10198.3329
$ws.Range("L126").Value = 250015440
$ws.Range("M126").Value = -7728.332900000001
$ws.Range("N126").Value = -250020380

$ws.Range("H132").Value = 2797.125
$ws.Range("I132").Value = 2005.28
$ws.Range("K132").Value = 6015.84
$ws.Range("M132").Value = -3485.84

$ws.Range("H134").Value = 1259.4667
$ws.Range("I134").Value = 932.44446
$ws.Range("J134").Value = 1750
$ws.Range("K134").Value = 2797.33338
$ws.Range("L134").Value = 5250
$ws.Range("M134").Value = -262.33338
$ws.Range("N134").Value = -10320

$ws.Range("H136").Value = 20273.482
$ws.Range("I136").Value = 1861.6
$ws.Range("J136").Value = 31104
$ws.Range("K136").Value = 5584.799999999999
$ws.Range("L136").Value = 93312
$ws.Range("M136").Value = -3034.799999999999
$ws.Range("N136").Value = -98412


# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1192.3889
$ws.Range("I5").Value = 754.0769
$ws.Range("J5").Value = 2332
$ws.Range("K5").Value = 2262.2307
$ws.Range("L5").Value = 6996
$ws.Range("M5").Value = -2150.2307
$ws.Range("N5").Value = -7220

$ws.Range("H37").Value = 250000000
$ws.Range("J37").Value = 250000000
$ws.Range("L37").Value = 750000000
$ws.Range("N37").Value = -750000224

$ws.Range("H131").Value = 102886.58
$ws.Range("I131").Value = 524.9
$ws.Range("J131").Value = 114518.59
$ws.Range("K131").Value = 1574.7
$ws.Range("L131").Value = 343555.77
$ws.Range("M131").Value = 3465.3
$ws.Range("N131").Value = -353635.77

$ws.Range("H135").Value = 1192.3889
$ws.Range("I135").Value = 754.0769
$ws.Range("J135").Value = 2332
$ws.Range("K135").Value = 6786.6921
$ws.Range("L135").Value = 20988
$ws.Range("M135").Value = -4251.6921
$ws.Range("N135").Value = -26058

$ws.Range("H140").Value = 1562.6842
$ws.Range("J140").Value = 2999.8572
$ws.Range("L140").Value = 8999.571599999999
$ws.Range("N140").Value = -19359.5716


# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1610.9412
$ws.Range("I97").Value = 1785.4166
$ws.Range("J97").Value = 1192.2
$ws.Range("K97").Value = 1785.4166
$ws.Range("L97").Value = 1192.2
$ws.Range("M97").Value = -1289.4166
$ws.Range("N97").Value = -2184.2

$ws.Range("H126").Value = 4307
$ws.Range("I126").Value = 3675.6667
$ws.Range("J126").Value = 4938.3335
$ws.Range("K126").Value = 11027.0001
$ws.Range("L126").Value = 14815.0005
$ws.Range("M126").Value = -8557.000100000001
$ws.Range("N126").Value = -19755.0005

$ws.Range("H132").Value = 45052.668
$ws.Range("I132").Value = 4200
$ws.Range("J132").Value = 65479
$ws.Range("K132").Value = 12600
$ws.Range("L132").Value = 196437
$ws.Range("M132").Value = -10070
$ws.Range("N132").Value = -201497


# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2483.2778
$ws.Range("I68").Value = 2344.5557
$ws.Range("J68").Value = 2622
$ws.Range("K68").Value = 2344.5557
$ws.Range("L68").Value = 2622
$ws.Range("M68").Value = -1595.5557
$ws.Range("N68").Value = -4120

$ws.Range("H71").Value = 2483.2778
$ws.Range("I71").Value = 2344.5557
$ws.Range("J71").Value = 2622
$ws.Range("K71").Value = 11722.7785
$ws.Range("L71").Value = 13110
$ws.Range("M71").Value = -7978.7785
$ws.Range("N71").Value = -20598

$ws.Range("H132").Value = 1341838.2
$ws.Range("I132").Value = 2010424.4
$ws.Range("J132").Value = 4666
$ws.Range("K132").Value = 6031273.199999999
$ws.Range("L132").Value = 13998
$ws.Range("M132").Value = -6028743.199999999
$ws.Range("N132").Value = -19058

$ws.Range("H136").Value = 1908.1818
$ws.Range("I136").Value = 1665
$ws.Range("K136").Value = 4995
$ws.Range("M136").Value = -2445


# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3709.0908
$ws.Range("I62").Value = 2260
$ws.Range("J62").Value = 4916.6665
$ws.Range("K62").Value = 2260
$ws.Range("L62").Value = 4916.6665
$ws.Range("M62").Value = -1636
$ws.Range("N62").Value = -6164.6665

$ws.Range("H65").Value = 3709.0908
$ws.Range("I65").Value = 2260
$ws.Range("J65").Value = 4916.6665
$ws.Range("K65").Value = 11300
$ws.Range("L65").Value = 24583.3325
$ws.Range("M65").Value = -8180
$ws.Range("N65").Value = -30823.3325

$ws.Range("H81").Value = 76924190
$ws.Range("I81").Value = 1312.5
$ws.Range("K81").Value = 2625
$ws.Range("M81").Value = -1564

$ws.Range("H84").Value = 76924190
$ws.Range("I84").Value = 1312.5
$ws.Range("K84").Value = 13125
$ws.Range("M84").Value = -7821

$ws.Range("H107").Value = 6494992.5
$ws.Range("I107").Value = 1733.8334
$ws.Range("J107").Value = 45454544
$ws.Range("K107").Value = 5201.5002
$ws.Range("L107").Value = 136363632
$ws.Range("M107").Value = -3281.5002
$ws.Range("N107").Value = -136367472

$ws.Range("H132").Value = 1658.3334
$ws.Range("I132").Value = 907.6667
$ws.Range("K132").Value = 2723.0001
$ws.Range("M132").Value = -193.0001000000002

$ws.Range("H136").Value = 23461908
$ws.Range("I136").Value = 28674972
$ws.Range("K136").Value = 86024916
$ws.Range("M136").Value = -86022366

